$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44425
$ws.Range("J2").Value = 25
$ws.Range("L2").Value = 25000
$ws.Range("M2").Value = 24520
$ws.Range("P2").Value = 1635

# Row 3
$ws.Range("D3").Value = 44446
$ws.Range("J3").Value = 34

# Row 4
$ws.Range("D4").Value = 44349
$ws.Range("J4").Value = 21
$ws.Range("K4").Value = 24000
$ws.Range("L4").Value = 25000
$ws.Range("M4").Value = 24524
$ws.Range("P4").Value = 1635

# Row 5
$ws.Range("D5").Value = 44677
$ws.Range("J5").Value = 34
$ws.Range("K5").Value = 25000
$ws.Range("L5").Value = 26000
$ws.Range("M5").Value = 25500
$ws.Range("P5").Value = 1700

# Row 6
$ws.Range("D6").Value = 44460
$ws.Range("J6").Value = 25
$ws.Range("M6").Value = 24480
$ws.Range("P6").Value = 1632

# Row 7
$ws.Range("D7").Value = 44719
$ws.Range("J7").Value = 43
$ws.Range("K7").Value = 17000
$ws.Range("L7").Value = 18000
$ws.Range("M7").Value = 17512
$ws.Range("P7").Value = 1167

# Row 8
$ws.Range("D8").Value = 44329
$ws.Range("J8").Value = 25
$ws.Range("K8").Value = 23000
$ws.Range("L8").Value = 23000
$ws.Range("M8").Value = 23000
$ws.Range("P8").Value = 1533

# Row 9
$ws.Range("D9").Value = 44383
$ws.Range("J9").Value = 25
$ws.Range("K9").Value = 13000
$ws.Range("L9").Value = 14000
$ws.Range("M9").Value = 13480
$ws.Range("P9").Value = 899

# Row 10
$ws.Range("D10").Value = 44397
$ws.Range("K10").Value = 23000
$ws.Range("L10").Value = 24000
$ws.Range("M10").Value = 23500
$ws.Range("P10").Value = 1567

# Row 11
$ws.Range("D11").Value = 44428
$ws.Range("J11").Value = 16
$ws.Range("K11").Value = 25000
$ws.Range("M11").Value = 25500
$ws.Range("P11").Value = 1700

# Row 12
$ws.Range("D12").Value = 44680
$ws.Range("J12").Value = 36
$ws.Range("K12").Value = 24000
$ws.Range("L12").Value = 25000
$ws.Range("M12").Value = 24500
$ws.Range("P12").Value = 1633

# Row 13
$ws.Range("D13").Value = 44406
$ws.Range("J13").Value = 25
$ws.Range("K13").Value = 24000
$ws.Range("L13").Value = 25000
$ws.Range("M13").Value = 24520
$ws.Range("P13").Value = 1635

# Row 14
$ws.Range("D14").Value = 44413

# Row 15
$ws.Range("D15").Value = 44729
$ws.Range("J15").Value = 52
$ws.Range("L15").Value = 24000
$ws.Range("M15").Value = 24000
$ws.Range("P15").Value = 1600

# Row 16
$ws.Range("D16").Value = 44341
$ws.Range("J16").Value = 36
$ws.Range("K16").Value = 24000
$ws.Range("L16").Value = 25000
$ws.Range("M16").Value = 24500
$ws.Range("P16").Value = 1633

# Row 17
$ws.Range("D17").Value = 44727
$ws.Range("J17").Value = 28
$ws.Range("L17").Value = 24000
$ws.Range("M17").Value = 24000
$ws.Range("P17").Value = 1600

# Row 18
$ws.Range("D18").Value = 44708
$ws.Range("J18").Value = 25
$ws.Range("K18").Value = 26000
$ws.Range("L18").Value = 26000
$ws.Range("M18").Value = 26000
$ws.Range("P18").Value = 1733

# Row 19
$ws.Range("D19").Value = 44705
$ws.Range("J19").Value = 35
$ws.Range("K19").Value = 26000
$ws.Range("L19").Value = 26000
$ws.Range("M19").Value = 26000
$ws.Range("P19").Value = 1733

# Row 20
$ws.Range("D20").Value = 44400
$ws.Range("J20").Value = 16
$ws.Range("K20").Value = 24000
$ws.Range("L20").Value = 25000
$ws.Range("M20").Value = 24500
$ws.Range("P20").Value = 1633

# Row 21
$ws.Range("D21").Value = 44343
$ws.Range("J21").Value = 26
$ws.Range("K21").Value = 23000
$ws.Range("M21").Value = 23500
$ws.Range("P21").Value = 1567

# Row 22
$ws.Range("D22").Value = 44453
$ws.Range("J22").Value = 25
$ws.Range("K22").Value = 25000
$ws.Range("L22").Value = 26000
$ws.Range("M22").Value = 25520
$ws.Range("P22").Value = 1701

# Row 23
$ws.Range("D23").Value = 44432

# Row 24
$ws.Range("D24").Value = 44463
$ws.Range("J24").Value = 25
$ws.Range("M24").Value = 24480
$ws.Range("P24").Value = 1632

# Row 25
$ws.Range("D25").Value = 44714
$ws.Range("J25").Value = 52
$ws.Range("K25").Value = 18000
$ws.Range("L25").Value = 20000
$ws.Range("M25").Value = 19000
$ws.Range("P25").Value = 1267

# Row 26
$ws.Range("D26").Value = 44351
$ws.Range("J26").Value = 34
$ws.Range("K26").Value = 24000
$ws.Range("L26").Value = 25000
$ws.Range("M26").Value = 24500
$ws.Range("P26").Value = 1633

# Row 27
$ws.Range("D27").Value = 44685
$ws.Range("J27").Value = 20
$ws.Range("K27").Value = 25000
$ws.Range("L27").Value = 25000
$ws.Range("M27").Value = 25000
$ws.Range("P27").Value = 1667

# Row 28
$ws.Range("D28").Value = 44449
$ws.Range("J28").Value = 18
$ws.Range("M28").Value = 24500
$ws.Range("P28").Value = 1633

# Row 29
$ws.Range("D29").Value = 44455
$ws.Range("J29").Value = 18
$ws.Range("L29").Value = 25000
$ws.Range("M29").Value = 24500
$ws.Range("P29").Value = 1633

# Row 30
$ws.Range("D30").Value = 44435
$ws.Range("J30").Value = 34
$ws.Range("K30").Value = 24000
$ws.Range("L30").Value = 25000
$ws.Range("M30").Value = 24500
$ws.Range("P30").Value = 1633

# Row 31
$ws.Range("D31").Value = 44390
$ws.Range("J31").Value = 34
$ws.Range("K31").Value = 24000
$ws.Range("L31").Value = 25000
$ws.Range("M31").Value = 24500
$ws.Range("P31").Value = 1633

# Row 32
$ws.Range("D32").Value = 44706
$ws.Range("J32").Value = 30
$ws.Range("K32").Value = 26000
$ws.Range("L32").Value = 26000
$ws.Range("M32").Value = 26000
$ws.Range("P32").Value = 1733

# Row 33
$ws.Range("D33").Value = 44411
$ws.Range("K33").Value = 25000
$ws.Range("L33").Value = 26000
$ws.Range("M33").Value = 25500
$ws.Range("P33").Value = 1700

# Row 34
$ws.Range("D34").Value = 44418
$ws.Range("J34").Value = 16
$ws.Range("K34").Value = 25000
$ws.Range("L34").Value = 26000
$ws.Range("M34").Value = 25500
$ws.Range("P34").Value = 1700

# Row 35
$ws.Range("D35").Value = 44726
$ws.Range("J35").Value = 28
$ws.Range("K35").Value = 24000
$ws.Range("L35").Value = 24000
$ws.Range("M35").Value = 24000
$ws.Range("P35").Value = 1600

# Row 36
$ws.Range("D36").Value = 44336
$ws.Range("J36").Value = 34
$ws.Range("K36").Value = 24000
$ws.Range("L36").Value = 25000
$ws.Range("M36").Value = 24500
$ws.Range("P36").Value = 1633

# Row 37
$ws.Range("D37").Value = 44385
$ws.Range("K37").Value = 14000
$ws.Range("L37").Value = 15000
$ws.Range("M37").Value = 14480
$ws.Range("P37").Value = 965

# Row 38
$ws.Range("D38").Value = 44707
$ws.Range("J38").Value = 30
$ws.Range("K38").Value = 26000
$ws.Range("L38").Value = 26000
$ws.Range("M38").Value = 26000
$ws.Range("P38").Value = 1733

# Row 39
$ws.Range("D39").Value = 44421
$ws.Range("J39").Value = 18
$ws.Range("K39").Value = 24000
$ws.Range("L39").Value = 25000
$ws.Range("M39").Value = 24500
$ws.Range("P39").Value = 1633

# Row 40
$ws.Range("D40").Value = 44442
$ws.Range("J40").Value = 28
$ws.Range("K40").Value = 24000
$ws.Range("L40").Value = 25000
$ws.Range("M40").Value = 24500
$ws.Range("P40").Value = 1633
